$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (and shared string header) from 08-12 to 08-13
$ws.Name = "Through 2022-08-13"
$ws.Range("I1").Value = "2022 (through 08-13)"

# Update August 2022 value (I9) and the Total row (I14)
$ws.Range("I9").Value = 75
$ws.Range("I14").Value = 1045
